{"js": "// Replace each three-digit \u00f7 one-digit division answer with its new value.\n// [oldText, newText] pairs, applied in document order via Body.search + Range.insertText(\"Replace\").\nconst replacements = [\n  [\"842\u00f76=140, 2\", \"123\u00f79=13, 6\"],\n  [\"144\u00f76=24, 0\", \"508\u00f77=72, 4\"],\n  [\"346\u00f74=86, 2\", \"159\u00f73=53, 0\"],\n  [\"748\u00f75=149, 3\", \"329\u00f78=41, 1\"],\n  [\"542\u00f76=90, 2\", \"638\u00f72=319, 0\"],\n  [\"140\u00f73=46, 2\", \"887\u00f79=98, 5\"],\n  [\"983\u00f78=122, 7\", \"425\u00f78=53, 1\"],\n  [\"151\u00f76=25, 1\", \"853\u00f76=142, 1\"],\n  [\"922\u00f73=307, 1\", \"881\u00f79=97, 8\"],\n  [\"890\u00f74=222, 2\", \"597\u00f79=66, 3\"],\n  [\"535\u00f79=59, 4\", \"686\u00f78=85, 6\"],\n  [\"699\u00f78=87, 3\", \"731\u00f76=121, 5\"],\n  [\"506\u00f74=126, 2\", \"385\u00f76=64, 1\"],\n  [\"545\u00f76=90, 5\", \"794\u00f78=99, 2\"],\n  [\"932\u00f72=466, 0\", \"515\u00f79=57, 2\"],\n  [\"270\u00f73=90, 0\", \"142\u00f74=35, 2\"],\n  [\"896\u00f79=99, 5\", \"342\u00f74=85, 2\"],\n  [\"232\u00f72=116, 0\", \"747\u00f77=106, 5\"],\n  [\"544\u00f74=136, 0\", \"882\u00f73=294, 0\"],\n  [\"360\u00f75=72, 0\", \"673\u00f79=74, 7\"],\n  [\"980\u00f78=122, 4\", \"468\u00f74=117, 0\"],\n  [\"839\u00f79=93, 2\", \"656\u00f73=218, 2\"],\n  [\"982\u00f79=109, 1\", \"454\u00f77=64, 6\"],\n  [\"497\u00f74=124, 1\", \"263\u00f78=32, 7\"],\n  [\"276\u00f79=30, 6\", \"432\u00f79=48, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  found.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit \u00f7 one-digit division answer with its new value.\n# @(oldText, newText) pairs, applied via Range.Find.Execute(..., Replace:=wdReplaceAll).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"842\u00f76=140, 2\", \"123\u00f79=13, 6\"),\n    @(\"144\u00f76=24, 0\", \"508\u00f77=72, 4\"),\n    @(\"346\u00f74=86, 2\", \"159\u00f73=53, 0\"),\n    @(\"748\u00f75=149, 3\", \"329\u00f78=41, 1\"),\n    @(\"542\u00f76=90, 2\", \"638\u00f72=319, 0\"),\n    @(\"140\u00f73=46, 2\", \"887\u00f79=98, 5\"),\n    @(\"983\u00f78=122, 7\", \"425\u00f78=53, 1\"),\n    @(\"151\u00f76=25, 1\", \"853\u00f76=142, 1\"),\n    @(\"922\u00f73=307, 1\", \"881\u00f79=97, 8\"),\n    @(\"890\u00f74=222, 2\", \"597\u00f79=66, 3\"),\n    @(\"535\u00f79=59, 4\", \"686\u00f78=85, 6\"),\n    @(\"699\u00f78=87, 3\", \"731\u00f76=121, 5\"),\n    @(\"506\u00f74=126, 2\", \"385\u00f76=64, 1\"),\n    @(\"545\u00f76=90, 5\", \"794\u00f78=99, 2\"),\n    @(\"932\u00f72=466, 0\", \"515\u00f79=57, 2\"),\n    @(\"270\u00f73=90, 0\", \"142\u00f74=35, 2\"),\n    @(\"896\u00f79=99, 5\", \"342\u00f74=85, 2\"),\n    @(\"232\u00f72=116, 0\", \"747\u00f77=106, 5\"),\n    @(\"544\u00f74=136, 0\", \"882\u00f73=294, 0\"),\n    @(\"360\u00f75=72, 0\", \"673\u00f79=74, 7\"),\n    @(\"980\u00f78=122, 4\", \"468\u00f74=117, 0\"),\n    @(\"839\u00f79=93, 2\", \"656\u00f73=218, 2\"),\n    @(\"982\u00f79=109, 1\", \"454\u00f77=64, 6\"),\n    @(\"497\u00f74=124, 1\", \"263\u00f78=32, 7\"),\n    @(\"276\u00f79=30, 6\", \"432\u00f79=48, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    # Fresh range over the whole document each time so every search starts from the top.\n    $range = $d.Content\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap:=wdFindContinue(1), Format, ReplaceWith, Replace:=wdReplaceAll(2))\n    $result = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n"}
